$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 (header): add two new columns, OPK2 / OPK3 ----
$ws.Cells.Item(1, 3).Value = "OPK2"
$ws.Cells.Item(1, 4).Value = "OPK3"

# ---- Row 2: "MO" becomes "Wydział MO" ----
$ws.Cells.Item(2, 2).Value = "Wydział MO"

# Cells A3, A4, A5, A6, A7, A8 hold values that look numeric ("73.54", "10", ...).
# Pre-formatting each as Text before assigning keeps them as real text (same
# as typing into a Text-formatted cell in Excel); ClearFormats() afterwards
# removes the temporary formatting again so no extra style is left applied.

# ---- Row 3 ----
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "73.54"
$ws.Cells.Item(3, 1).ClearFormats()
$ws.Cells.Item(3, 2).Value = "Wydział MO"
$ws.Cells.Item(3, 3).Value = "Linia MO"
$ws.Cells.Item(3, 4).Value = "Stanowisko 54"

# ---- Row 4 ----
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "73.55"
$ws.Cells.Item(4, 1).ClearFormats()
$ws.Cells.Item(4, 2).Value = "Wydział MO"
$ws.Cells.Item(4, 3).Value = "Linia MO"
$ws.Cells.Item(4, 4).Value = "Stanowisko 55"

# ---- Row 5 ----
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "73.516"
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(5, 2).Value = "Wydział MO"
$ws.Cells.Item(5, 3).Value = "Linia MO"
$ws.Cells.Item(5, 4).Value = "Stanowisko 516"

# ---- Row 6 ----
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "73.515"
$ws.Cells.Item(6, 1).ClearFormats()
$ws.Cells.Item(6, 2).Value = "Wydział MO"
$ws.Cells.Item(6, 3).Value = "Linia MO"
$ws.Cells.Item(6, 4).Value = "Stanowisko 515"

# ---- Row 7 ----
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "10"
$ws.Cells.Item(7, 1).ClearFormats()
$ws.Cells.Item(7, 2).Value = "Wydział LAK"
$ws.Cells.Item(7, 3).Value = "Hala LAK"

# ---- Row 8 ----
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "10.1"
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 2).Value = "Wydział LAK"
$ws.Cells.Item(8, 3).Value = "Hala LAK"

# ---- Row 9 ---- ("10.1.1" is not numeric-parsable, so it stays text naturally)
$ws.Cells.Item(9, 1).Value = "10.1.1"
$ws.Cells.Item(9, 2).Value = "Wydział LAK"
$ws.Cells.Item(9, 3).Value = "Hala LAK"
$ws.Cells.Item(9, 4).Value = "Kabina A"
